$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Interface")

# --- Add new column F ("JurisdictionWiseReport") ---

# Copy formatting from column E into column F for the used rows (header + data)
$ws.Range("E1:E11").Copy()
$ws.Range("F1:F11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Header
$ws.Range("F1").Value = "JurisdictionWiseReport"

# Data rows: default to "No" for the new column
$ws.Range("F2").Value = "No"
$ws.Range("F3").Value = "No"
$ws.Range("F4").Value = "No"
$ws.Range("F5").Value = "No"
$ws.Range("F6").Value = "No"
$ws.Range("F7").Value = "No"
$ws.Range("F8").Value = "No"
$ws.Range("F9").Value = "No"
$ws.Range("F10").Value = "No"
$ws.Range("F11").Value = "Yes"

# Row 9 Execute flag toggled on (previously "No")
$ws.Range("E9").Value = "Yes"
# Row 11 Execute flag toggled off (previously "Yes"), since the Agent Module
# execute flag moved into the new JurisdictionWiseReport column
$ws.Range("E11").Value = "No"

# Update the conditional formatting range to include the new column
$fc = $ws.Range("D2:E11").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("D2:F11"))

# Update selection to reflect the saved cursor position
$ws.Range("J10").Select()
